# "contingencies with rene fine" - extend the parallel-case grid by two more
# columns (P, Q) and rebalance the I/K/M/O contingency columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header row (row 1): append P1=14, Q1=15, matching the existing
#     bold/bordered/centered header style used by B1:O1. ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial($xlPasteFormats)

# --- Data rows (rows 2-25): the I/K/M/O columns swap their 1/2 values,
#     and two new columns P and Q are appended with value 2. ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column
}
